$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" sheet and
# the "全部类型" sheet, which carry duplicate rows for the same events:
#   F8  (苏州·漫语堂动漫嘉年华): 253 -> 254
#   F11 (国泰北路18号 GTC卡丁车场 event): 16 -> 17

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F8").Value = 254
    $ws.Range("F11").Value = 17
}
